# Update gh-pages to output generated at 456a3b4
# Sheet "展览" (index 1) and sheet "全部类型" (index 4) both contain the same
# underlying event rows; the "F" column ("想去人数") counters were refreshed.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # 展览
$ws1.Range("F2").Value = 4338
$ws1.Range("F4").Value = 335
$ws1.Range("F7").Value = 44
$ws1.Range("F9").Value = 130
$ws1.Range("F10").Value = 312
$ws1.Range("F11").Value = 245
$ws1.Range("F12").Value = 2948
$ws1.Range("F13").Value = 145
$ws1.Range("F14").Value = 1519
$ws1.Range("F15").Value = 10

$ws4 = $wb.Worksheets.Item(4)   # 全部类型
$ws4.Range("F2").Value = 4338
$ws4.Range("F4").Value = 335
$ws4.Range("F8").Value = 44
$ws4.Range("F10").Value = 130
$ws4.Range("F11").Value = 312
$ws4.Range("F12").Value = 245
$ws4.Range("F13").Value = 2948
$ws4.Range("F14").Value = 145
$ws4.Range("F15").Value = 1519
$ws4.Range("F16").Value = 10
